$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.482.36'
$ws.Range("D3").Value = '3.152.92'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '612.73'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '143.99'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.149.36'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").Value = '5.35'
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").Value = '35.53'
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").Value = '3.673.44'
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '64.457.92'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '3.154.02'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '6.85'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = '476.65'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '14.67'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '0.724'
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '84.67'
$ws.Range("E25").Value = '  +2.05%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("D28").Value = '8.55'
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").Value = '7.37'
$ws.Range("E29").Value = '  +7.79%  '
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  +2.36%  '
$ws.Range("E31").Value = '  -4.35%  '
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("D33").Value = '26.46'
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  +1.21%  '
$ws.Range("D36").Value = '5.95'
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").Value = '52.77'
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("D38").Value = '0.0₃0749'
$ws.Range("E38").Value = '  +3.32%  '
$ws.Range("D39").Value = '3.11'
$ws.Range("E39").Value = '  +4.44%  '
$ws.Range("D40").Value = '452.55'
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").Value = '0.0396'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").Value = '2.845.36'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("E47").Value = '  +6.30%  '
$ws.Range("D48").Value = '26.44'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '120.87'
$ws.Range("E51").Value = '  +1.93%  '
